# AutomationTest-FinalProject: Write Test Cases For Employee Reports Function
# Adds a new "EmployeeReports" worksheet (copied from "EditCustomer" so it
# inherits the same look & feel / conditional styling), then fills it in
# with the new Employee Report test cases.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Create the new sheet by copying "EditCustomer" (keeps fonts/fills/
#    borders/column widths close to the existing "sister" sheets) and
#    drop it in after the last sheet, then rename it.
# ---------------------------------------------------------------------
$template = $wb.Worksheets.Item("EditCustomer")
$template.Copy([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws.Name = "EmployeeReports"

# ---------------------------------------------------------------------
# 2. Wipe the old (copied) data so we start from a clean A1:K12 block,
#    then widen the used range out to K.
# ---------------------------------------------------------------------
$ws.Range("A1:K12").Clear()

# ---------------------------------------------------------------------
# 3. Column widths (user resized every column on this sheet).
# ---------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 22.27
$ws.Columns.Item(2).ColumnWidth = 34.54
$ws.Columns.Item(3).ColumnWidth = 22.27
$ws.Columns.Item(4).ColumnWidth = 22.27
$ws.Columns.Item(5).ColumnWidth = 28.54
$ws.Columns.Item(6).ColumnWidth = 21.18
$ws.Columns.Item(7).ColumnWidth = 21.0
$ws.Columns.Item(8).ColumnWidth = 17.73
$ws.Columns.Item(9).ColumnWidth = 33.09
$ws.Columns.Item(10).ColumnWidth = 33.73
$ws.Columns.Item(11).ColumnWidth = 22.27

# ---------------------------------------------------------------------
# 4. Header row.
# ---------------------------------------------------------------------
$ws.Range("A1").Value = "TCs"
$ws.Range("B1").Value = "Description"
$ws.Range("C1").Value = "Employee Name"
$ws.Range("D1").Value = "Project Name"
$ws.Range("E1").Value = "Activity Name"
$ws.Range("F1").Value = "Start Date"
$ws.Range("G1").Value = "End Date"
$ws.Range("H1").Value = "Approved Timesheets"
$ws.Range("I1").Value = "Expected"
$ws.Range("J1").Value = "Actual"
$ws.Range("K1").Value = "Result"
$ws.Range("A1:K1").Font.Bold = $true
$ws.Range("A1:K1").HorizontalAlignment = -4108
$ws.Range("A1:K1").VerticalAlignment = -4108
$ws.Range("A1:K1").WrapText = $true
$ws.Rows.Item(1).RowHeight = 33

# ---------------------------------------------------------------------
# 5. Test-case rows (2-10).
# ---------------------------------------------------------------------
$tcs = @("TC_OHRM_EP_01","TC_OHRM_EP_02","TC_OHRM_EP_03","TC_OHRM_EP_04","TC_OHRM_EP_05","TC_OHRM_EP_06","TC_OHRM_EP_07","TC_OHRM_EP_08","TC_OHRM_EP_09")

$descriptions = @(
  "Verify that system displays the Employee Report after the user selected an employee.",
  "Verify that system displays the Employee Report after the user selected an employee and the project name.",
  "Verify that system displays the Employee Report after the user selected an employee, the project name, and a specific activity in the project.",
  "Verify that system displays the Employee Report after the user selected an employee, the project name, and the project date range in the project.",
  "Verify that system displays the Employee Report after the user selected an employee, the project name, a specific activity, and project date range in the project.",
  "Verify that system displays the Employee Report after the user selected an employee, and turn on 'Only Include Approved Timesheet'.",
  "Verify that system displays the Employee Report after the user selected an employee, the project name, and turns on 'Only Include Approved Timesheet'.",
  "Verify that system displays the Employee Report after the user selected an employee, the project name, project's activity, and turns on 'Only Include Approved Timesheet'.",
  "Verify that system displays the Employee Report after the user selected an employee, the project name, the project's activity, project date range, and turns on 'Only Include Approved Timesheet'."
)

$employee = "a"
$project = "ACME"
$activity = "Administration"
$startDate = "2023-08-13"
$endDate = "2023-08-20"
$expectedActual = "orangehrm-paper-container"

# Which rows (1-based within 2..10) have project / activity / date-range / approved-flag filled in.
$hasProject  = @(0,1,1,1,1,0,1,1,1)
$hasActivity = @(0,0,1,0,1,0,0,1,1)
$hasDates    = @(0,0,0,1,1,0,0,0,1)
$hasApproved = @(0,0,0,0,0,1,1,1,1)

for ($i = 0; $i -lt 9; $i++) {
    $r = $i + 2
    $ws.Rows.Item($r).RowHeight = 154.5

    $ws.Range("A$r").Value = $tcs[$i]
    $ws.Range("B$r").Value = $descriptions[$i]

    $ws.Range("C$r").NumberFormat = "@"
    $ws.Range("C$r").Value = $employee

    if ($hasProject[$i] -eq 1) {
        $ws.Range("D$r").Value = $project
    }
    if ($hasActivity[$i] -eq 1) {
        $ws.Range("E$r").Value = $activity
    }
    if ($hasDates[$i] -eq 1) {
        $ws.Range("F$r").NumberFormat = "@"
        $ws.Range("F$r").Value = $startDate
        $ws.Range("G$r").NumberFormat = "@"
        $ws.Range("G$r").Value = $endDate
    }
    if ($hasApproved[$i] -eq 1) {
        $ws.Range("H$r").Formula = "=LOWER(TRUE)"
    }

    $ws.Range("I$r").Value = $expectedActual
    $ws.Range("J$r").Value = $expectedActual
    $ws.Range("K$r").Value = "PASSED"

    $ws.Range("A$r`:K$r").HorizontalAlignment = -4108
    $ws.Range("A$r`:K$r").VerticalAlignment = -4108
    $ws.Range("A$r`:K$r").WrapText = $true
}

# ---------------------------------------------------------------------
# 6. Summary row + trailing blank row (matches the template's
#    "TCs Passed / Total TCs" style footer row).
# ---------------------------------------------------------------------
$ws.Rows.Item(11).RowHeight = 154.5
$ws.Range("H11").Value = "9/10"
$ws.Range("H11").HorizontalAlignment = -4108
$ws.Range("H11").VerticalAlignment = -4108
$ws.Range("H11").WrapText = $true

$ws.Rows.Item(12).RowHeight = 154.5

# ---------------------------------------------------------------------
# 7. View state: zoom out a bit, select J2, make this the active sheet.
# ---------------------------------------------------------------------
$ws.Activate()
$win = $excel.ActiveWindow
$win.Zoom = 55
$ws.Range("J2").Select()
